# Wallet ledger update - "data till 4 Mar 3AM"
# Adds two new ledger rows (85, 86) continuing the running-balance chain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 85: 2020-03-03 (serial 43893), an "Ordered Amount" debit of 38480 ---
$ws.Range("A85").Value = 43893
$ws.Range("A85").NumberFormat = '[$-409]d\-mmm\-yyyy;@'
$ws.Range("B85").Value = 38480
$ws.Range("D85").Value = "Ordered Amount"
$ws.Range("E85").Formula = '=IF(A85="","",SUM(E84-B85+C85))'

# --- Row 86: 2020-03-04 (serial 43894), a "Manual Added" credit of 103900 ---
$ws.Range("A86").Value = 43894
$ws.Range("A86").NumberFormat = '[$-409]d\-mmm\-yyyy;@'
$ws.Range("C86").Value = 103900
$ws.Range("D86").Value = "Manual Added"
$ws.Range("E86").Formula = '=IF(A86="","",SUM(E85-B86+C86))'

# Move the frozen-pane selection down to where the new rows were entered.
$ws.Range("D88").Select()
